$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below the current (only) employee/period detail row (row 16),
# pushing the signature rows (old 21/22) down to 22/23.
$ws.Rows(17).Insert()

# Clone the formatting of row 16 (the detail row) into the newly inserted row 17.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New row keeps the same worker but the previous period (2506), while the
# existing row now reflects the newest period (2507).
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73114020"
$ws.Range("D17").Value = "GERARDO CARDONA ESPINOSA"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 877803

# The row that used to be the only period now shows the newest period.
$ws.Range("E16").Value = "2507"

# Update the summary fields: two periods now owed, so the total mora doubles.
$ws.Range("F13").Value = 2
$ws.Range("E11").Value = 113880
